$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 47671804
$ws.Range("I76").Value = 115285
$ws.Range("J76").Value = 83339190
$ws.Range("K76").Value = 115285
$ws.Range("L76").Value = 83339190
$ws.Range("M76").Value = -114970
$ws.Range("N76").Value = -83339820
$ws.Range("H79").Value = 47671804
$ws.Range("I79").Value = 115285
$ws.Range("J79").Value = 83339190
$ws.Range("K79").Value = 115285
$ws.Range("L79").Value = 83339190
$ws.Range("M79").Value = -114193
$ws.Range("N79").Value = -83341374
$ws.Range("H129").Value = 2129.6667
$ws.Range("I129").Value = 1024.5
$ws.Range("K129").Value = 3073.5
$ws.Range("M129").Value = 1926.5
$ws.Range("H138").Value = 4770.316
$ws.Range("I138").Value = 2115.8635
$ws.Range("K138").Value = 6347.5905
$ws.Range("M138").Value = -1207.5905

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I4").Value = 233.33333
$ws.Range("K4").Value = 233.33333
$ws.Range("M4").Value = -117.33333
$ws.Range("H6").Value = 24500
$ws.Range("I6").Value = 24500
$ws.Range("K6").Value = 24500
$ws.Range("M6").Value = -24327
$ws.Range("H132").Value = 198087.95
$ws.Range("I132").Value = 272389.8
$ws.Range("K132").Value = 817169.3999999999
$ws.Range("M132").Value = -814639.3999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2303.2307
$ws.Range("I105").Value = 2345.1667
$ws.Range("J105").Value = 1800
$ws.Range("K105").Value = 2345.1667
$ws.Range("L105").Value = 1800
$ws.Range("M105").Value = -598.1667000000002
$ws.Range("N105").Value = -5294
$ws.Range("H107").Value = 2932.3333
$ws.Range("I107").Value = 2818.8
$ws.Range("K107").Value = 2818.8
$ws.Range("M107").Value = -898.8000000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 491.21054
$ws.Range("I7").Value = 561.4375
$ws.Range("K7").Value = 561.4375
$ws.Range("M7").Value = -448.4375
$ws.Range("H31").Value = 695942.6
$ws.Range("J31").Value = 23183.79
$ws.Range("L31").Value = 23183.79
$ws.Range("N31").Value = -23773.79
$ws.Range("H34").Value = 695942.6
$ws.Range("J34").Value = 23183.79
$ws.Range("L34").Value = 23183.79
$ws.Range("N34").Value = -23587.79
$ws.Range("H53").Value = 70684
$ws.Range("J53").Value = 70684
$ws.Range("L53").Value = 70684
$ws.Range("N53").Value = -71898
$ws.Range("H99").Value = 6624.9165
$ws.Range("I99").Value = 4933.1665
$ws.Range("K99").Value = 4933.1665
$ws.Range("M99").Value = -3435.1665
$ws.Range("H126").Value = 6624.9165
$ws.Range("I126").Value = 4933.1665
$ws.Range("K126").Value = 14799.4995
$ws.Range("M126").Value = -12329.4995
$ws.Range("H132").Value = 5410.778
$ws.Range("I132").Value = 10200
$ws.Range("K132").Value = 30600
$ws.Range("M132").Value = -28070
$ws.Range("H134").Value = 424564.4
$ws.Range("I134").Value = 3036.3333
$ws.Range("K134").Value = 9108.999899999999
$ws.Range("M134").Value = -6573.999899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18582114
$ws.Range("I4").Value = 2648367.8
$ws.Range("K4").Value = 7945103.399999999
$ws.Range("M4").Value = -7944991.399999999
$ws.Range("H12").Value = 175.9375
$ws.Range("J12").Value = 220.08333
$ws.Range("L12").Value = 660.24999
$ws.Range("N12").Value = -1006.24999
$ws.Range("H92").Value = 792.4
$ws.Range("I92").Value = 279.5
$ws.Range("J92").Value = 1134.3334
$ws.Range("K92").Value = 838.5
$ws.Range("L92").Value = 3403.0002
$ws.Range("M92").Value = 409.5
$ws.Range("N92").Value = -5899.0002
$ws.Range("H112").Value = 146640.72
$ws.Range("I112").Value = 202977.2
$ws.Range("K112").Value = 608931.6000000001
$ws.Range("M112").Value = -607823.6000000001
$ws.Range("H137").Value = 2265.4707
$ws.Range("I137").Value = 1731.7273
$ws.Range("J137").Value = 3244
$ws.Range("K137").Value = 5195.1819
$ws.Range("L137").Value = 9732
$ws.Range("M137").Value = -95.18189999999959
$ws.Range("N137").Value = -19932

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 511229.9
$ws.Range("I7").Value = 839216.2
$ws.Range("J7").Value = 19250.5
$ws.Range("K7").Value = 839216.2
$ws.Range("L7").Value = 19250.5
$ws.Range("M7").Value = -839104.2
$ws.Range("N7").Value = -19474.5
$ws.Range("H16").Value = 2367.9
$ws.Range("J16").Value = 1999
$ws.Range("L16").Value = 1999
$ws.Range("N16").Value = -2339
$ws.Range("H46").Value = 2725.3057
$ws.Range("I46").Value = 2278.3635
$ws.Range("J46").Value = 3427.6428
$ws.Range("K46").Value = 2278.3635
$ws.Range("L46").Value = 3427.6428
$ws.Range("M46").Value = -2090.3635
$ws.Range("N46").Value = -3803.6428
$ws.Range("H74").Value = 36500
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 36500
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H100").Value = 68144.64999999999
$ws.Range("I100").Value = 77568.5
$ws.Range("J100").Value = 24166.666
$ws.Range("K100").Value = 77568.5
$ws.Range("L100").Value = 24166.666
$ws.Range("M100").Value = -77027.5
$ws.Range("N100").Value = -25248.666
$ws.Range("H122").Value = 791398.9
$ws.Range("I122").Value = 4430
$ws.Range("K122").Value = 13290
$ws.Range("M122").Value = -10840
$ws.Range("H126").Value = 511229.9
$ws.Range("I126").Value = 839216.2
$ws.Range("J126").Value = 19250.5
$ws.Range("K126").Value = 2517648.6
$ws.Range("L126").Value = 57751.5
$ws.Range("M126").Value = -2515178.6
$ws.Range("N126").Value = -62691.5
$ws.Range("H132").Value = 6801.4546
$ws.Range("I132").Value = 5438.0713
$ws.Range("J132").Value = 9187.375
$ws.Range("K132").Value = 16314.2139
$ws.Range("L132").Value = 27562.125
$ws.Range("M132").Value = -13784.2139
$ws.Range("N132").Value = -32622.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 7150
$ws.Range("J33").Value = 7150
$ws.Range("L33").Value = 7150
$ws.Range("N33").Value = -7650
$ws.Range("H36").Value = 7150
$ws.Range("J36").Value = 7150
$ws.Range("L36").Value = 7150
$ws.Range("N36").Value = -7650
$ws.Range("H130").Value = 73306.336
$ws.Range("J130").Value = 73306.336
$ws.Range("L130").Value = 73306.336
$ws.Range("N130").Value = -83346.336
$ws.Range("H132").Value = 28472.705
$ws.Range("I132").Value = 1895.5385
$ws.Range("K132").Value = 5686.6155
$ws.Range("M132").Value = -3156.6155
$ws.Range("H136").Value = 359494.88
$ws.Range("I136").Value = 373860.12
$ws.Range("K136").Value = 1121580.36
$ws.Range("M136").Value = -1119030.36
